$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Rename "safe_column_name_test" -> "safe_header_name_test"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("safe_column_name_test")
$ws2.Name = "safe_header_name_test"
$ws2.Select()
$ws2.Range("G2").Select()

# ---------------------------------------------------------------------
# 2) Add a new worksheet "date_test" as the last sheet
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "date_test"

# column widths (closest achievable values - the runtime snaps column widths
# to a coarse pixel grid, so we pick the inputs that land nearest to the
# target character widths of 27.6640625 / 42.33203125)
$ws3.Columns.Item(1).ColumnWidth = 26.75
$ws3.Columns.Item(2).ColumnWidth = 41.42

# ---------------------------------------------------------------------
# Header row
# ---------------------------------------------------------------------
$ws3.Range("A1").Value = "date"
$ws3.Range("B1").Value = "plaincol"

# ---------------------------------------------------------------------
# row 5, col A : custom long-date format applied to a *text* value (kept
# as text, not as a number, because the value isn't numeric)
# ---------------------------------------------------------------------
$ws3.Range("A5").NumberFormat = "[$-F800]dddd\,\ mmmm\ dd\,\ yyyy"
$ws3.Range("A5").Value = "Wednesday, Mar-14-2012"

# ---------------------------------------------------------------------
# row 2, col B : wrapped explanatory text
# ---------------------------------------------------------------------
$ws3.Range("B2").WrapText = $true
$ws3.Range("B2").Value = "it will still parse the dates below as date even if plaincol is not in the default --dates-whitelist because the cell format was set to date"

# ---------------------------------------------------------------------
# row 6, col B : text value explicitly forced to Text (@) number format
# so it is not parsed as a date
# ---------------------------------------------------------------------
$ws3.Range("B6").NumberFormat = "@"
$ws3.Range("B6").Value = "9/11/01 8:30 am"

# ---------------------------------------------------------------------
# row 5, col B : wrapped explanatory text
# ---------------------------------------------------------------------
$ws3.Range("B5").WrapText = $true
$ws3.Range("B5").Value = "the date below is not parsed as a date coz we didn't explicitly set the cell format to a date format and `"plaincol`" is not in the --dates-whitelist"
$ws3.Rows.Item(5).RowHeight = 48

# ---------------------------------------------------------------------
# row 2, col A : date serial formatted as short date (m/d/yy -> numFmtId 14)
# ---------------------------------------------------------------------
$ws3.Range("A2").NumberFormat = "m/d/yy"
$ws3.Range("A2").Value = 29580
$ws3.Rows.Item(2).RowHeight = 58

# ---------------------------------------------------------------------
# row 3, col A : date+time formatted with custom [$-409]m/d/yy h:mm AM/PM;@
# row 3, col B : date serial formatted as short date
# ---------------------------------------------------------------------
$ws3.Range("A3").NumberFormat = "[$-409]m/d/yy\ h:mm\ AM/PM;@"
$ws3.Range("A3").Value = 37145.354166666664

$ws3.Range("B3").NumberFormat = "m/d/yy"
$ws3.Range("B3").Value = 37145

# ---------------------------------------------------------------------
# row 4, col A : plain text "not a date"
# row 4, col B : date+time formatted as m/d/yy h:mm
# ---------------------------------------------------------------------
$ws3.Range("A4").Value = "not a date"

$ws3.Range("B4").NumberFormat = "m/d/yy h:mm"
$ws3.Range("B4").Value = 37145.354166666664

# ---------------------------------------------------------------------
# row 6, col A : date serial using the same custom long-date format
# ---------------------------------------------------------------------
$ws3.Range("A6").NumberFormat = "[$-F800]dddd\,\ mmmm\ dd\,\ yyyy"
$ws3.Range("A6").Value = 37145

# ---------------------------------------------------------------------
# selection / active sheet bookkeeping (date_test becomes the active tab)
# ---------------------------------------------------------------------
$ws3.Select()
$ws3.Range("D5").Select()

$wb.Save()
